$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.47
$ws.Range("D2").Value = -0.19
$ws.Range("E2").Value = 0.59
$ws.Range("F2").Value = 0.39

$ws.Range("C3").Value = -0.01
$ws.Range("D3").Value = 0.01
$ws.Range("E3").Value = -0.07
$ws.Range("F3").Value = 0.42

$ws.Range("C4").Value = 0.31
$ws.Range("D4").Value = -0.01
$ws.Range("E4").Value = 0.25
$ws.Range("F4").Value = 0.18

$ws.Range("C5").Value = 0.55
$ws.Range("D5").Value = 0.27
$ws.Range("E5").Value = 0.68
$ws.Range("F5").Value = 0.49

$ws.Range("C6").Value = 0.36
$ws.Range("D6").Value = -0.15
$ws.Range("E6").Value = 0.21
$ws.Range("F6").Value = 0.03

$ws.Range("C7").Value = 0.97
$ws.Range("D7").Value = 0.09
$ws.Range("E7").Value = 0.72
$ws.Range("F7").Value = 0.52

$ws.Range("C8").Value = 0.29
$ws.Range("D8").Value = 0.22
$ws.Range("E8").Value = 0.1
$ws.Range("F8").Value = -0.02

$ws.Range("C9").Value = 0.06
$ws.Range("D9").Value = 0.25
$ws.Range("E9").Value = 0.22
$ws.Range("F9").Value = -0.05

$ws.Range("C10").Value = 0.22
$ws.Range("D10").Value = 0.04
$ws.Range("E10").Value = -0.26
$ws.Range("F10").Value = 0.24

$ws.Range("C11").Value = 0.54
$ws.Range("D11").Value = 0.36
$ws.Range("E11").Value = 0.08
$ws.Range("F11").Value = 0.09

$ws.Range("C12").Value = -0.72
$ws.Range("D12").Value = 0.23
$ws.Range("E12").Value = -0.19
$ws.Range("F12").Value = 0.03

$ws.Range("C13").Value = 0.13
$ws.Range("D13").Value = 0.68
$ws.Range("E13").Value = -0.16
$ws.Range("F13").Value = -0.15

$ws.Range("C14").Select()
